$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ethnicities")

# Add a "Category" header label in A1 (previously empty)
$ws.Range("A1").Value = "Category"

# A new data column is being inserted as the "F" column (between the
# existing "Other" and "Unknown" ethnicity columns). The values that used
# to live in F2:F5 move over to the newly appended G2:G5 cells, and the
# new column's figures are written into F2:F5.
$ws.Range("G2").Value = $ws.Range("F2").Value()
$ws.Range("G3").Value = $ws.Range("F3").Value()
$ws.Range("G4").Value = $ws.Range("F4").Value()
$ws.Range("G5").Value = $ws.Range("F5").Value()

$ws.Range("F2").Value = 92
$ws.Range("F3").Value = 17
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 2

# Update the active selection on the sheet
$ws.Range("E10").Select()

Write-Host "Edit applied"
